$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 24, shifting existing rows 24:93 down to 25:94.
$ws.Rows(24).Insert()

# Populate the newly inserted row 24 with the new data record.
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = 45044
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 100112043
$ws.Range("G24").Value = "Pepino dulce"
$ws.Range("H24").Value = "Cultivar IV Región"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 20000
$ws.Range("N24").Value = '$/bandeja 18 kilos'
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 1111
$ws.Range("Q24").Value = 18
$ws.Range("R24").Value = "Hortaliza"
